$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49 (shifts rows 49-95 down to 50-96,
# and copies formatting -- e.g. the date style in column D -- from the
# surrounding rows).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly data record.
$ws.Range("A49").Value = 6
$ws.Range("B49").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44781
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = 100114007
$ws.Range("G49").Value = "Jengibre"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 260
$ws.Range("K49").Value = 1300
$ws.Range("L49").Value = 12000
$ws.Range("M49").Value = 6238
$ws.Range("N49").Value = "$/caja 13 kilos"
$ws.Range("O49").Value = "Perú"
$ws.Range("P49").Value = 480
$ws.Range("Q49").Value = 13
$ws.Range("R49").Value = "Hortaliza"
